# Feat: implement component Table in readEXCEL
# Rename the "Colegio" table column to "Institucion".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Editing the header cell directly renames the table column ("Colegio" ->
# "Institucion") the same way Excel does when you type over a table header.
$ws.Range("D1").Value = "Institucion"

# Move the active selection to D2, matching the edited workbook state.
$ws.Range("D2").Select()
